# Update the correction-factor spreadsheet:
#  - Replace the hard-coded pi approximation (3.1415) with PI() in the
#    circumference/correction-factor formulas in column H.
#  - Add a new "Uncertainty" calculation column I that propagates the
#    relative uncertainties of B/C, D/E and F/G through H.
#  - Leave the final selection on I2 (matches the saved workbook state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H: replace the hard-coded pi (3.1415) with PI().
$ws.Range("H2").Formula = "=F2/((B2/20)^2*PI()*D2/10)"
$ws.Range("H3:H4").Formula = "=F3/((B3/20)^2*PI()*D3/10)"
$ws.Range("H7:H8").Formula = "=F7/((B7/20)^2*PI()*D7/10)"

# Column I: new uncertainty-propagation calculation. Rows 5 and 6 have no
# foil circumference/weight data (B/C are blank there), so that formula
# group spans I3:I8 but rows 5-6 are cleared back out, leaving only
# I2 (standalone) and I3, I4, I7, I8 (populated).
$ws.Range("I2").Formula = "=SQRT((C2/B2)^2+(E2/D2)^2+(G2/F2)^2)*H2"
$ws.Range("I3:I8").Formula = "=SQRT((C3/B3)^2+(E3/D3)^2+(G3/F3)^2)*H3"
$ws.Range("I5").ClearContents()
$ws.Range("I6").ClearContents()

$ws.Range("I2").Select()
